# The MoxR family ATPase was added to the list of queries.
#
# Net effect on the deck: the statistics table on slide 1 was rebuilt
# (PowerPoint hands a freshly (re)created/duplicated graphicFrame a new
# shape id + default "Table N" name) and slightly enlarged/repositioned,
# two cell values were incremented to reflect the new query, and the row
# heights/column widths were re-flowed by PowerPoint's layout engine to
# their new autofit extents.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$orig = $s.Shapes.Item(1)

# Duplicate the existing table so the new shape keeps every bit of the
# original formatting (borders, fills, merges, fonts, ...), then drop the
# original. The duplicate is the first table PowerPoint (re)creates in
# this session, so it is auto-named/auto-numbered "Table 1" / id 2 --
# exactly like a freshly (re)inserted table would be.
$dup = $orig.Duplicate()
$orig.Delete()
$tbl_shape = $dup
$tbl_shape.Name = "Table 1"

$tbl = $tbl_shape.Table

# Reposition/resize the table (values below are the exact point-valued
# equivalents of the target EMU offsets/extents).
$tbl_shape.Left = 2.8333858267716536
$tbl_shape.Top = 2.6667718110236223

# Column widths (EMU 614815 / 1522848 / 1021537 x5).
$tbl.Columns.Item(1).Width = 48.41063118110236
$tbl.Columns.Item(2).Width = 119.90929133858268
$tbl.Columns.Item(3).Width = 80.43598559055118
$tbl.Columns.Item(4).Width = 80.43598559055118
$tbl.Columns.Item(5).Width = 80.43598559055118
$tbl.Columns.Item(6).Width = 80.43598559055118
$tbl.Columns.Item(7).Width = 80.43598559055118

# Row heights (EMU 638115 for the header row, 208364 for every other row).
$tbl.Rows.Item(1).Height = 50.24527559055118
for ($i = 2; $i -le $tbl.Rows.Count; $i++) {
    $tbl.Rows.Item($i).Height = 16.406614173228345
}

# New query (MoxR family ATPase) bumps the "Other" row and the TOTAL row
# in the "Number of medium subunit genes" column.
$tbl.Cell(12, 5).Shape.TextFrame.TextRange.Text = "47"
$tbl.Cell(13, 5).Shape.TextFrame.TextRange.Text = "1476"
